# Fixed sample portfolio 1 delta matrix
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# EUR Government Bonds delta (row 8, col B) was 1 -> 0.01
$ws.Range("B8").Value = 0.01

# EUR Investment Grade Corporate Bonds delta (row 10, col D) was 1 -> 0.01
$ws.Range("D10").Value = 0.01

# Update the active selection to match the authored state
$ws.Range("D11").Select()
